$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: mark as unselected -> strip yellow highlight style, clear the
# empty D2/G2 placeholder cells entirely, and set K2 (status) to 0.
$ws.Range("A2:C2").Style = "Normal"
$ws.Range("D2").Clear()
$ws.Range("E2:F2").Style = "Normal"
$ws.Range("G2").Clear()
$ws.Range("H2:J2").Style = "Normal"
$ws.Range("K2").Style = "Normal"
$ws.Range("K2").Value = 0

# Row 3: same treatment - unselect.
$ws.Range("A3:C3").Style = "Normal"
$ws.Range("D3").Clear()
$ws.Range("E3:F3").Style = "Normal"
$ws.Range("G3").Clear()
$ws.Range("H3:J3").Style = "Normal"
$ws.Range("K3").Style = "Normal"
$ws.Range("K3").Value = 0

# Row 6: mark as selected -> apply the yellow highlight style to the whole
# row (including previously-absent empty D6/G6 cells), and set K6 to 1.
$ws.Range("A6:C6").Interior.Color = 65535
$ws.Range("D6").Interior.Color = 65535
$ws.Range("E6:F6").Interior.Color = 65535
$ws.Range("G6").Interior.Color = 65535
$ws.Range("H6:J6").Interior.Color = 65535
$ws.Range("K6").Interior.Color = 65535
$ws.Range("K6").Value = 1
